# Reorder the tax-column headers on both data sheets so that the
# "Taxable Value / CGST / SGST / IGST" columns read
#   E=Taxable Value, F=IGST, G=CGST, H=SGST
# (previously F=CGST, G=SGST, H=IGST).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Item(1, 6).Value = "IGST"
    $ws.Cells.Item(1, 7).Value = "CGST"
    $ws.Cells.Item(1, 8).Value = "SGST"
}
